$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test01")

# --- Remove the second ("Sheet1") worksheet entirely ---
$wb.Worksheets.Item("Sheet1").Delete()

# --- Row 1 no longer carries an explicit custom row height ---
$ws.Rows.Item(1).AutoFit()

# --- A3 becomes a hyperlinked facebook URL cell, re-using the plain
#     "Hyperlink" cell format shared by the other linked cells (instead of
#     its old one-off Courier New style) ---
$ws.Hyperlinks.Add($ws.Range("A3"), "https://www.facebook.com/")
$ws.Range("A3").Value = "https://www.facebook.com/"
$ws.Range("B2").Copy()
$ws.Range("A3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- New row 4: facebook URL / hello (mailto test@gmail.com) / test@123456A ---
# Hyperlinks.Add calls happen in B4, C4, A4 order so the generated
# relationship ids (rId7, rId8, rId9) line up with the target workbook.
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:test@gmail.com", "", "", "test@gmail.com")
$ws.Range("B4").Value = "hello"

$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:test@123456A")
$ws.Range("C4").Value = "test@123456A"

$ws.Hyperlinks.Add($ws.Range("A4"), "https://www.facebook.com/")
$ws.Range("A4").Value = "https://www.facebook.com/"

$ws.Range("B2").Copy()
$ws.Range("B4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B2").Copy()
$ws.Range("C4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("B2").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Selection ends on B4 ---
$ws.Range("B4").Select()
